# ajout tableaux rdt et profils
$wb = $excel.ActiveWorkbook

# Rename the first sheet ("Feuille 1 - 1_continu_SS" -> "Continu Lent")
$wsLent = $wb.Worksheets.Item(1)
$wsLent.Name = "Continu Lent"

# The newly renamed sheet becomes the active/selected tab,
# with a new active cell selection at P18.
$wsLent.Activate()
$wsLent.Range("P18").Select()
